$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")

# Recapitalize the two lower-case entries used on the Elements sheet
$ws.Range("A2").Value = "Rating"
$ws.Range("A3").Value = "Contact Form"

# Move the cell selection to A3 (was B7)
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
